# Updates the cryptos list on the active worksheet with refreshed Price (D)
# and Volume(1h) (E) figures, and fixes the RenderToken / Decentraland row
# ordering (rows 44-45), matching the upstream "Updated cryptos list" commit.
#
# All target cells hold plain text in the source workbook (e.g. "29.154.11",
# "1.002", "  +0.08%  "). Most new values (multi-dot prices, percentages
# wrapped in spaces) are naturally kept as text by Excel. A handful of new
# Price values look like ordinary decimal numbers (e.g. "1.002", "325.56"),
# which Excel would otherwise auto-convert into numeric values, so for those
# cells we briefly force the cell's number format to Text ("@") before
# assigning the value, guaranteeing the text is stored exactly as in the
# source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param(
        [string]$CellRef,
        [string]$Text
    )
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
}

# --- Price / Volume(1h) refresh, rows 2-43 --------------------------------

$ws.Range("D2").Value = "29.154.11"
$ws.Range("E2").Value = "  +0.08%  "

$ws.Range("D3").Value = "1.904.81"
$ws.Range("E3").Value = "  +0.01%  "

Set-TextCell "D4" "1.002"
$ws.Range("E4").Value = "  -0.06%  "

Set-TextCell "D5" "325.56"
$ws.Range("E5").Value = "  -0.40%  "

Set-TextCell "D6" "1.001"
$ws.Range("E6").Value = "  -0.17%  "

Set-TextCell "D7" "0.4605"
$ws.Range("E7").Value = "  -0.22%  "

Set-TextCell "D8" "0.3893"
$ws.Range("E8").Value = "  -0.96%  "

Set-TextCell "D9" "0.07877"
$ws.Range("E9").Value = "  -0.72%  "

Set-TextCell "D10" "0.9903"
$ws.Range("E10").Value = "  -1.06%  "

Set-TextCell "D11" "21.93"
$ws.Range("E11").Value = "  -1.42%  "

$ws.Range("D12").Value = "1.883.67"
$ws.Range("E12").Value = "  +1.14%  "

Set-TextCell "D13" "5.777"
$ws.Range("E13").Value = "  +0.33%  "

Set-TextCell "D14" "7.054"
$ws.Range("E14").Value = "  -0.23%  "

Set-TextCell "D15" "0.07022"
$ws.Range("E15").Value = "  +1.18%  "

Set-TextCell "D16" "87.94"
$ws.Range("E16").Value = "  -0.44%  "

Set-TextCell "D17" "1.003"
$ws.Range("E17").Value = "  -0.03%  "

Set-TextCell "D18" "0.000009930"
$ws.Range("E18").Value = "  -1.27%  "

Set-TextCell "D19" "17.08"
$ws.Range("E19").Value = "  -0.20%  "

Set-TextCell "D20" "1.001"
$ws.Range("E20").Value = "  -0.12%  "

$ws.Range("D21").Value = "29.148.45"
$ws.Range("E21").Value = "  +0.04%  "

Set-TextCell "D22" "5.324"
$ws.Range("E22").Value = "  -0.60%  "

$ws.Range("E23").Value = "  +0.40%  "

Set-TextCell "D24" "2.101"
$ws.Range("E24").Value = "  +2.36%  "

Set-TextCell "D25" "155.97"
$ws.Range("E25").Value = "  -0.33%  "

Set-TextCell "D26" "19.41"
$ws.Range("E26").Value = "  -0.09%  "

Set-TextCell "D27" "5.898"
$ws.Range("E27").Value = "  -3.53%  "

Set-TextCell "D28" "118.44"
$ws.Range("E28").Value = "  -0.11%  "

$ws.Range("E29").Value = "  -6.16%  "

Set-TextCell "D30" "0.09327"
$ws.Range("E30").Value = "  -0.53%  "

Set-TextCell "D31" "0.8942"
$ws.Range("E31").Value = "  -3.68%  "

Set-TextCell "D32" "5.246"

$ws.Range("E33").Value = "  -2.02%  "

Set-TextCell "D34" "3.141"
$ws.Range("E34").Value = "  -3.95%  "

Set-TextCell "D35" "0.05790"
$ws.Range("E35").Value = "  -0.63%  "

Set-TextCell "D36" "1.167"
$ws.Range("E36").Value = "  -3.15%  "

Set-TextCell "D37" "0.02090"
$ws.Range("E37").Value = "  -0.67%  "

Set-TextCell "D38" "1.000"
$ws.Range("E38").Value = "  -0.11%  "

Set-TextCell "D39" "7.671"
$ws.Range("E39").Value = "  -3.06%  "

Set-TextCell "D40" "0.5682"
$ws.Range("E40").Value = "  -1.00%  "

Set-TextCell "D41" "0.1797"
$ws.Range("E41").Value = "  -0.04%  "

Set-TextCell "D42" "9.723"
$ws.Range("E42").Value = "  -2.22%  "

Set-TextCell "D43" "11.86"
$ws.Range("E43").Value = "  -0.54%  "

# --- Rows 44 & 45: RenderToken / Decentraland swap + refreshed figures ----
# Row 44 (rank 42) was RenderToken, now becomes Decentraland.
# Row 45 (rank 43) was Decentraland, now becomes RenderToken.

$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextCell "D44" "0.5355"
$ws.Range("E44").Value = "  -1.00%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D45" "2.198"
$ws.Range("E45").Value = "  -1.60%  "

# --- Remaining Price / Volume(1h) refresh, rows 46-51 ----------------------

Set-TextCell "D46" "0.07016"
$ws.Range("E46").Value = "  -0.78%  "

Set-TextCell "D47" "1.849"
$ws.Range("E47").Value = "  -1.49%  "

Set-TextCell "D48" "2.552"
$ws.Range("E48").Value = "  -0.03%  "

Set-TextCell "D49" "113.18"
$ws.Range("E49").Value = "  +0.08%  "

Set-TextCell "D50" "0.2941"
$ws.Range("E50").Value = "  +0.13%  "

$ws.Range("E51").Value = "  -2.71%  "

Write-Output "cryptos list updated"
